$wb = $excel.ActiveWorkbook

# Map of sheet name -> B2, C2, B3, C3 new values
$updates = @{
    "Study case 1" = @{ B2 = 0.75784908686924;  C2 = 3.928147304886926; B3 = 0.7452189676354882; C3 = 4.014194655748607 }
    "Study case 2" = @{ B2 = 0.7486132925533111; C2 = 4.646439359915873; B3 = 0.7341723396478861; C3 = 4.742345045752013 }
    "Study case 3" = @{ B2 = 0.75784908686924;  C2 = 3.928147304886926; B3 = 0.7274824667297396; C3 = 4.177093172467385 }
    "Study case 4" = @{ B2 = 0.75784908686924;  C2 = 3.928147304886926; B3 = 0.7274824667297396; C3 = 4.177093172467385 }
    "Experiment 1" = @{ B2 = 0.6624752232246374; C2 = 4.462123810427518; B3 = 0.6560882265792992; C3 = 4.493955455642348 }
    "Experiment 2" = @{ B2 = 0.6624752232246374; C2 = 4.462123810427518; B3 = 0.6472810699802006; C3 = 4.550112781227658 }
    "Experiment 3" = @{ B2 = 0.7070534090477629; C2 = 5.180111994054941; B3 = 0.6886481381166791; C3 = 5.316211754647637 }
    "Experiment 4" = @{ B2 = 0.7070534090477629; C2 = 5.180111994054941; B3 = 0.698739271031948;  C3 = 5.233174198540799 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $vals = $updates[$sheetName]
    $ws.Range("B2").Value = $vals.B2
    $ws.Range("C2").Value = $vals.C2
    $ws.Range("B3").Value = $vals.B3
    $ws.Range("C3").Value = $vals.C3
}
